$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header row: "Call Name *" -> "Call / Distribution Name *"
# ---------------------------------------------------------------
$ws.Range("C1").Value = "Call / Distribution Name *"

# ---------------------------------------------------------------
# Data rows. Fund renamed "Macro Fund" -> "SAAS Fund", Folio No
# replaced with small call/distribution numbers, and existing
# "Call 2" rows become "Call 1" / new "Distribution 1" rows are
# appended below as negative-quantity reversals.
# ---------------------------------------------------------------
$rows = @(
    @{ Row=2;  Fund="SAAS Fund"; Num=6;  Name="Call 1";         Unit="Series C"; Qty=10;  Price=100; Prem=10; Year=2022 },
    @{ Row=3;  Fund="SAAS Fund"; Num=7;  Name="Call 1";         Unit="Series C"; Qty=20;  Price=100; Prem=10; Year=2022 },
    @{ Row=4;  Fund="SAAS Fund"; Num=8;  Name="Call 1";         Unit="Series A"; Qty=30;  Price=100; Prem=0;  Year=2022 },
    @{ Row=5;  Fund="SAAS Fund"; Num=9;  Name="Call 1";         Unit="Series A"; Qty=40;  Price=100; Prem=0;  Year=2022 },
    @{ Row=6;  Fund="SAAS Fund"; Num=10; Name="Call 1";         Unit="Series B"; Qty=50;  Price=100; Prem=5;  Year=2022 },
    @{ Row=7;  Fund="SAAS Fund"; Num=6;  Name="Distribution 1"; Unit="Series C"; Qty=-5;  Price=100; Prem=10; Year=2023 },
    @{ Row=8;  Fund="SAAS Fund"; Num=7;  Name="Distribution 1"; Unit="Series C"; Qty=-10; Price=100; Prem=10; Year=2023 },
    @{ Row=9;  Fund="SAAS Fund"; Num=8;  Name="Distribution 1"; Unit="Series A"; Qty=-15; Price=100; Prem=0;  Year=2023 },
    @{ Row=10; Fund="SAAS Fund"; Num=9;  Name="Distribution 1"; Unit="Series A"; Qty=-20; Price=100; Prem=0;  Year=2023 },
    @{ Row=11; Fund="SAAS Fund"; Num=10; Name="Distribution 1"; Unit="Series B"; Qty=-25; Price=100; Prem=5;  Year=2023 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Fund
    $ws.Cells.Item($row, 1).Font.Name = "Arial"

    $ws.Cells.Item($row, 2).Value = $r.Num
    $ws.Cells.Item($row, 2).Font.Name = "Arial"
    $ws.Cells.Item($row, 2).Font.Size = 10

    $ws.Cells.Item($row, 3).Value = $r.Name
    $ws.Cells.Item($row, 4).Value = $r.Unit
    $ws.Cells.Item($row, 5).Value = $r.Qty
    $ws.Cells.Item($row, 6).Value = $r.Price
    $ws.Cells.Item($row, 7).Value = $r.Prem

    # Set the date format before the value so freshly-created rows reuse
    # the existing short-date style instead of inferring a new numFmt.
    $ws.Cells.Item($row, 8).NumberFormat = "m/d/yy"
    $ws.Cells.Item($row, 8).Value = (Get-Date -Year $r.Year -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0)

    $ws.Cells.Item($row, 9).Value = "No"
    $ws.Cells.Item($row, 10).Value = "Upload"
}

# ---------------------------------------------------------------
# Column widths (auto-fit changed slightly after the new content)
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10.1328125
$ws.Columns.Item(2).ColumnWidth = 8.73046875
$ws.Columns.Item(3).ColumnWidth = 11

# ---------------------------------------------------------------
# Selection ends up on the last entered cell
# ---------------------------------------------------------------
$ws.Range("H11").Select()
